$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# --- New row 5 data on the "ShareSkill" sheet ---
# Title
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Value = "Galaxy Test Inspector#"

# Description
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Value = "Software and Hardware Testing^"

# SkillExchange (set before Tag so shared-string order matches source file)
$ws.Range("P5").WrapText = $true
$ws.Range("P5").Value = "Jmeter"

# Tag
$ws.Range("H5").WrapText = $true
$ws.Range("H5").Value = "Performance Testing"

# ServiceType
$ws.Range("J5").WrapText = $true
$ws.Range("J5").Value = "Hourly"

# LocationType
$ws.Range("K5").WrapText = $true
$ws.Range("K5").Value = "Online"

# StartDate
$ws.Range("L5").WrapText = $true
$ws.Range("L5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("L5").Value = (Get-Date -Year 2022 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)

# EndDate
$ws.Range("M5").WrapText = $true
$ws.Range("M5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("M5").Value = (Get-Date -Year 2023 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)

# The new row wraps onto two lines, same as rows 2-4 above it
$ws.Rows.Item(5).RowHeight = 27.6

# Match the selection saved with the workbook
$ws.Range("A5").Select() | Out-Null
